# Participant-specific stimulation order workbook:
# - row 2 (block 1) of channels/electrodes was already validated/correct
# - rows 3-5 get their final, validated channel counts and electrode-pair
#   sequences, and the "test block" rows (4 and 5) are re-colored to match
#   the normal light-green row styling used elsewhere in the sheet
# - cursor position is left where the author last clicked

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (block 3): re-color to the light-green style used by the other rows ---
$ws.Range("B4:E4").Interior.Color = 15073253
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("E4").Value = "[(1, 2)]"
$ws.Range("D4").Value = "[1]"

# --- Row 5 (block 4): re-color to the light-green style too ---
$ws.Range("B5:E5").Interior.Color = 15073253
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = "[2]"
$ws.Range("E5").Value = "[(3, 4)]"

# --- Row 3 (block 2): channel count + new trial/electrode sequences ---
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = "[2, 1]"
$ws.Range("E3").Value = "[(3, 4), (1, 2)]"

# --- Move the active selection, matching the author's last cursor position ---
$ws.Range("D12").Select()
